$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new date column BD with header "23-ago" (follows existing "22-ago" in BC)
$ws.Range("BD1").Value = "23-ago"

# Fill the new column's data values for rows 2-18
$values = @{
    2  = 0
    3  = 10.462988970924169
    4  = 18.307513916545272
    5  = 22.329040240787812
    6  = 0
    7  = 14.647075100176433
    8  = 19.910911972361824
    9  = 12.996003441818548
    10 = 5.5697887286241388
    11 = 11.372584698734656
    12 = 0
    13 = 10.3897971017286
    14 = 0
    15 = 0
    16 = 5.0108105179688893
    17 = 0
    18 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("BD$row").Value = $values[$row]
}

# Update the active selection to match the saved view state
$ws.Range("BE6").Select()
